$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 09:14"

# Row 55 - Ghana
$ws.Range("B55").Value = 43505
$ws.Range("D55").Value = 41532
$ws.Range("E55").Value = 1712

# Row 57 - Armenia
$ws.Range("B57").Value = 42792
$ws.Range("C57").Value = 176
$ws.Range("D57").Value = 35991
$ws.Range("E57").Value = 5949
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 852

# Row 63 - Afganistan
$ws.Range("B63").Value = 37999
$ws.Range("C63").Value = 46
$ws.Range("D63").Value = 28180
$ws.Range("E63").Value = 8432
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 1387

# Row 72 - Australia
$ws.Range("B72").Value = 24812
$ws.Range("C72").Value = 210
$ws.Range("E72").Value = 5551

# Row 73 - El Salvador
$ws.Range("D73").Value = 12032
$ws.Range("E73").Value = 11727
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 661

# Row 151 - Letonia
$ws.Range("B151").Value = 1337
$ws.Range("C151").Value = 4
$ws.Range("E151").Value = 211
